# Updates to power sector for curtailment and moving CHP out of flexible
# resources; updates to fuel balancing priorities.
#
# On the "FPIEBP" sheet, swap the priority-1 / priority-3 columns (B / C)
# for: petroleum gasoline (row 10), petroleum diesel (row 11),
# jet fuel or kerosene (row 14), heavy fuel oil (row 19) and
# LPG propane or butane (row 20).

$wb = $excel.ActiveWorkbook
$wsFPIEBP = $wb.Worksheets.Item("FPIEBP")

$rows = @(10, 11, 14, 19, 20)
foreach ($r in $rows) {
    $bCell = $wsFPIEBP.Cells.Item($r, 2)
    $cCell = $wsFPIEBP.Cells.Item($r, 3)
    $bVal = $bCell.Value2
    $cVal = $cCell.Value2
    $bCell.Value = $cVal
    $cCell.Value = $bVal
}

# Make FPIEBP the active sheet/tab, with B20:D20 selected (matches the
# selection left behind by the edit), and leave "About" unselected.
$wsFPIEBP.Activate()
$wsFPIEBP.Range("B20:D20").Select()
